# Auto-generated PowerShell script to apply the FICHAMENTO.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Step 1: copy cell formatting (style) from row 169 (same A,B,G,H,I layout) to new rows 171-183 ---
$templateRow = 169
$newRows = 171..183
foreach ($r in $newRows) {
    foreach ($col in @("A","B","G","H","I")) {
        $ws.Range("$col$templateRow").Copy()
        $ws.Range("$col$r").PasteSpecial(-4122) | Out-Null
    }
}
$excel.CutCopyMode = 0

# --- Step 2: set cell values in the precise chronological order that produces the same new-shared-string order as the original edit ---
$ws.Range("I171").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"
$ws.Range("A171").Value2 = "Scenario building: Uses and abuses"
$ws.Range("B171").Value2 = "A scenario is not a future reality but rather a means to represent it with the aim of clarifying present action in light of possible and desirable futures."
$ws.Range("B172").Value2 = "Thus, without a careful and attentive reader, many scenarios pass as credible as if the reader is guilty of not having understood the underlying meaning."
$ws.Range("B173").Value2 = "A scenario is a description (usually of a possble future) which assumes the intervention of several key events or conditions which will have taken place between the time of the original situation and the time in which the scwenario is set."
$ws.Range("B174").Value2 = "A scenario must satisfy the following conditions: pertinence, coherence, likelihood, importance and transparency."
$ws.Range("B175").Value2 = "Exploratory scenarios are concerned with past and present trends and lead to likely futures., Normative scenarios are constructed from alternative images of the future which may be both desirable and feare, and are conceived in a retro projective way. Thus, exploratory scenarios are devoid of human values, ewhereas normative scenarios are the expression of human values."
$ws.Range("B176").Value2 = "A scenario is not and end in itself. I only has meaning as na aid to decision-making in so far as it clarifies the consequences of current decisions."
$ws.Range("B177").Value2 = "Scenario splanning requires time to be done right, and a 12- to 18-month timeframe is not rare. Time is needed in prospective and strategic diagnoses so that scenarios involving the environment may be developed and that the main stakes for a specific organization based on possible futures, as revealed by scenarios, may be reviewed."
$ws.Range("B180").Value2 = "As a result, it is important to gather as many informed judgments as possible and then forge a consensus."
$ws.Range("B181").Value2 = "The uncertainty of the future can be evaluated across a number of scenarios which share he field of probable futures."
$ws.Range("B182").Value2 = "However, it is imporatnt to take into consideration the content of the various scenarios since the more probable among them may be either very similar or quite contrasted to one another. In theory, two possible situations may present themselves."
$ws.Range("B183").Value2 = "Using morphological analysis, a global system can be decomposed into dimensions (key questions concerning the future). Theses dimensions are: demographis, economic, technological, and social/organizational."
$ws.Range("H177").Value2 = "Prazo"
$ws.Range("H178").Value2 = "Prazo"
$ws.Range("B178").Value2 = "In extreme cases, policy-makers may launch a Foresight study that they wish to see finished in a matter of weeks. In this event, the prevailing conditions are rarely ideal, thought it is better to light a candle than curse the darkness."
$ws.Range("B179").Value2 = "Given a short time-frame, it is often advisable to limit the scenarios to several key hypotheses, say four to six."

# --- Step 3: fill remaining cells that reuse already-existing shared strings (order does not matter for these) ---
$ws.Range("G171").Value2 = "Cenários"
$ws.Range("H171").Value2 = "Conceito"
$ws.Range("A172").Value2 = "Scenario building: Uses and abuses"
$ws.Range("G172").Value2 = "Cenários"
$ws.Range("H172").Value2 = "Caracteristica"
$ws.Range("I172").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"
$ws.Range("A173").Value2 = "Scenario building: Uses and abuses"
$ws.Range("G173").Value2 = "Cenários "
$ws.Range("H173").Value2 = "Conceito"
$ws.Range("I173").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"
$ws.Range("A174").Value2 = "Scenario building: Uses and abuses"
$ws.Range("G174").Value2 = "Cenários"
$ws.Range("H174").Value2 = "Caracteristica"
$ws.Range("I174").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"
$ws.Range("A175").Value2 = "Scenario building: Uses and abuses"
$ws.Range("G175").Value2 = "Cenários"
$ws.Range("H175").Value2 = "Tipos"
$ws.Range("I175").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"
$ws.Range("A176").Value2 = "Scenario building: Uses and abuses"
$ws.Range("G176").Value2 = "Cenários"
$ws.Range("H176").Value2 = "Caracteristica"
$ws.Range("I176").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"
$ws.Range("A177").Value2 = "Scenario building: Uses and abuses"
$ws.Range("G177").Value2 = "Cenários"
$ws.Range("I177").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"
$ws.Range("A178").Value2 = "Scenario building: Uses and abuses"
$ws.Range("G178").Value2 = "Cenários"
$ws.Range("I178").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"
$ws.Range("A179").Value2 = "Scenario building: Uses and abuses"
$ws.Range("G179").Value2 = "Cenários"
$ws.Range("H179").Value2 = "Quantidade de cenários"
$ws.Range("I179").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"
$ws.Range("A180").Value2 = "Scenario building: Uses and abuses"
$ws.Range("G180").Value2 = "Cenários"
$ws.Range("H180").Value2 = "Quantidade de cenários"
$ws.Range("I180").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"
$ws.Range("A181").Value2 = "Scenario building: Uses and abuses"
$ws.Range("G181").Value2 = "Cenários"
$ws.Range("H181").Value2 = "Quantidade de cenários"
$ws.Range("I181").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"
$ws.Range("A182").Value2 = "Scenario building: Uses and abuses"
$ws.Range("G182").Value2 = "Cenários"
$ws.Range("H182").Value2 = "Quantidade de cenários"
$ws.Range("I182").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"
$ws.Range("A183").Value2 = "Scenario building: Uses and abuses"
$ws.Range("G183").Value2 = "Cenários"
$ws.Range("H183").Value2 = "Tipos"
$ws.Range("I183").Value2 = "@article{article,`nauthor = {Durance, Philippe and Godet, Michel},`nyear = {2010},`nmonth = {11},`npages = {1488-1492},`ntitle = {Scenario building: Uses and abuses},`nvolume = {77},`njournal = {Technological Forecasting and Social Change - TECHNOL FORECAST SOC CHANGE},`ndoi = {10.1016/j.techfore.2010.06.007}`n}"

# --- Step 4: set row heights for the new rows ---
$ws.Rows.Item(171).RowHeight = 28
$ws.Rows.Item(172).RowHeight = 28
$ws.Rows.Item(173).RowHeight = 42
$ws.Rows.Item(174).RowHeight = 28
$ws.Rows.Item(175).RowHeight = 70
$ws.Rows.Item(176).RowHeight = 28
$ws.Rows.Item(177).RowHeight = 56
$ws.Rows.Item(178).RowHeight = 42
$ws.Rows.Item(179).RowHeight = 28
$ws.Rows.Item(180).RowHeight = 28
$ws.Rows.Item(181).RowHeight = 28
$ws.Rows.Item(182).RowHeight = 42
$ws.Rows.Item(183).RowHeight = 42

# --- Step 5: set row heights for pre-existing rows whose auto-height changed (Excel version rendering difference) ---
$ws.Rows.Item(6).RowHeight = 98
$ws.Rows.Item(36).RowHeight = 84
$ws.Rows.Item(94).RowHeight = 168
$ws.Rows.Item(127).RowHeight = 154
$ws.Rows.Item(139).RowHeight = 42
$ws.Rows.Item(170).RowHeight = 84

# --- Step 6: update the view state (selection + scroll position) to match the target ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 168
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("I187").Select()

